# Weekly update: a new price record for "Puerro" (week of 2021-11-24,
# serial 44524) is inserted at the top of the data block (row 64),
# pushing all the existing rows (old 64..79) down by one (new 65..80).
#
# This mirrors how the upstream consolidation script prepends the newest
# week's observation and keeps the historical rows below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 64; everything from 64..79 shifts to 65..80.
$ws.Rows.Item(64).Insert()

# Populate the new row 64 with the latest week's record.
$ws.Cells.Item(64, 1).Value2  = 9
$ws.Cells.Item(64, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(64, 3).Value2  = "Metropolitana"
$ws.Cells.Item(64, 4).Value2  = 44524
$ws.Cells.Item(64, 5).Value2  = 13
$ws.Cells.Item(64, 6).Value2  = 100112005
$ws.Cells.Item(64, 7).Value2  = "Puerro"
$ws.Cells.Item(64, 8).Value2  = "Sin especificar"
$ws.Cells.Item(64, 9).Value2  = "Primera"
$ws.Cells.Item(64, 10).Value2 = 160
$ws.Cells.Item(64, 11).Value2 = 6000
$ws.Cells.Item(64, 12).Value2 = 7000
$ws.Cells.Item(64, 13).Value2 = 6500
$ws.Cells.Item(64, 14).Value2 = "`$/paquete 20 unidades"
$ws.Cells.Item(64, 15).Value2 = "Provincia de Chacabuco"
$ws.Cells.Item(64, 16).Value2 = 325
$ws.Cells.Item(64, 17).Value2 = 20
$ws.Cells.Item(64, 18).Value2 = "Hortaliza"
